$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Alignment: left + wrap text, applied to the whole sheet first so the
# shared/default style (used by every cell, including still-blank ones)
# picks up the new alignment.
$ws.Cells.HorizontalAlignment = -4131  # xlLeft
$ws.Cells.WrapText = $true

# New header row
$ws.Range("A1").Value = "Team List"
$ws.Range("B1").Value = "Location List"

# Column A: team letters A..K (rows 2-12)
$letters = @("A","B","C","D","E","F","G","H","I","J","K")
for ($i = 0; $i -lt $letters.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $letters[$i]
}

# Column B: location numbers 1..4 (rows 2-5 only)
for ($i = 1; $i -le 4; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $i
}

# Clear the previous "Games:" helper column content (column B rows 6-12) that
# is not part of the new layout. Clear() (not just ClearContents()) so the
# cells are removed entirely rather than left behind as empty/styled cells.
$ws.Range("B6:B12").Clear()

# Selection matches the saved state in the target file (A12)
$ws.Range("A12").Select()
